# Updated symbol list on Mon Jan 16 20:29:01 UTC 2023 with GitHub Actions
# Refresh the Price (column D) and Volume(1h) (column E) snapshot values
# for the coin rows on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (matches the sheet's existing
# inline-string cells for Price/Volume, e.g. "299.90", "-0.29%") without
# leaving the cell tagged as Text-formatted/quote-prefixed afterwards.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "299.90"
Set-TextValue $ws.Range("E2") "-0.29%"
Set-TextValue $ws.Range("D3") "31.78"
Set-TextValue $ws.Range("E3") "1.37%"
Set-TextValue $ws.Range("D4") "5.124"
Set-TextValue $ws.Range("E4") "0.19%"
Set-TextValue $ws.Range("D5") "0.08208"
Set-TextValue $ws.Range("E5") "11.46%"
Set-TextValue $ws.Range("D6") "2.573"
Set-TextValue $ws.Range("E6") "6.43%"
Set-TextValue $ws.Range("D7") "7.864"
Set-TextValue $ws.Range("E7") "-1.05%"
Set-TextValue $ws.Range("D8") "3.854"
Set-TextValue $ws.Range("E8") "1.72%"
Set-TextValue $ws.Range("E9") "1.19%"
Set-TextValue $ws.Range("D10") "0.1759"
Set-TextValue $ws.Range("E10") "3.03%"
Set-TextValue $ws.Range("D11") "0.07513"
Set-TextValue $ws.Range("E11") "-0.34%"
Set-TextValue $ws.Range("D12") "0.09051"
Set-TextValue $ws.Range("E12") "12.16%"
Set-TextValue $ws.Range("D13") "0.03016"
Set-TextValue $ws.Range("E13") "0.10%"
Set-TextValue $ws.Range("D14") "0.1001"
Set-TextValue $ws.Range("E14") "0.92%"
Set-TextValue $ws.Range("D15") "0.001525"
Set-TextValue $ws.Range("E15") "2.05%"
Set-TextValue $ws.Range("D16") "0.005936"
Set-TextValue $ws.Range("E16") "-2.84%"
Set-TextValue $ws.Range("D17") "3.621"
Set-TextValue $ws.Range("E17") "4.47%"
Set-TextValue $ws.Range("E19") "-1.20%"
Set-TextValue $ws.Range("D20") "0.1347"
Set-TextValue $ws.Range("E20") "0.81%"
Set-TextValue $ws.Range("D21") "4.134"
Set-TextValue $ws.Range("E21") "-11.02%"
Set-TextValue $ws.Range("D22") "0.1678"
Set-TextValue $ws.Range("E22") "7.13%"
Set-TextValue $ws.Range("D23") "0.04630"
Set-TextValue $ws.Range("E23") "-0.38%"
Set-TextValue $ws.Range("D24") "0.001247"
Set-TextValue $ws.Range("E24") "1.69%"
Set-TextValue $ws.Range("D25") "0.004548"
Set-TextValue $ws.Range("E25") "1.44%"
Set-TextValue $ws.Range("E26") "-7.60%"
Set-TextValue $ws.Range("D27") "0.0003404"
Set-TextValue $ws.Range("E27") "81.85%"
Set-TextValue $ws.Range("D39") "0.01783"
Set-TextValue $ws.Range("E39") "3.48%"
Set-TextValue $ws.Range("D40") "0.04592"
Set-TextValue $ws.Range("E40") "1.79%"
Set-TextValue $ws.Range("D41") "0.006918"
Set-TextValue $ws.Range("E41") "-4.49%"
Set-TextValue $ws.Range("E42") "2.58%"
Set-TextValue $ws.Range("D43") "0.002140"
Set-TextValue $ws.Range("E43") "-3.91%"
Set-TextValue $ws.Range("D44") "0.009846"
Set-TextValue $ws.Range("E44") "-8.15%"
Set-TextValue $ws.Range("D45") "0.00006184"
Set-TextValue $ws.Range("E45") "-1.71%"
Set-TextValue $ws.Range("E46") "-0.04%"
Set-TextValue $ws.Range("D47") "0.8059"
Set-TextValue $ws.Range("E47") "-57.83%"
Set-TextValue $ws.Range("D48") "0.008384"
Set-TextValue $ws.Range("E48") "-16.19%"
Set-TextValue $ws.Range("D49") "0.00002099"
Set-TextValue $ws.Range("E49") "-0.04%"
Set-TextValue $ws.Range("D50") "0.0001999"
Set-TextValue $ws.Range("E50") "0.03%"
